$wb = $excel.ActiveWorkbook

# --- Sheet "BigInt" (sheet1): fix mislabeled Return type column ---
# Rows 5-16 list BigInt operator overloads whose "Return type" was
# mistakenly written as "BigNum" instead of "BigInt".
$ws1 = $wb.Worksheets.Item("BigInt")
$ws1.Range("B5:B16").Value = "BigInt"

# Restore the default (unbolded/"Normal") cell style on B6, B8, B9 - in the
# original file they carried an explicit style that differed from their
# neighbours for no visual reason; the fix normalizes them back to the
# same (default) style as the rest of the column.
$ws1.Range("B6").Font.Bold = $false
$ws1.Range("B8").Font.Bold = $false
$ws1.Range("B9").Font.Bold = $false

# Reviewer reassignment: E15 and E16 move from "Bảo" to "Chiến"
$ws1.Range("E15").Value = "Chiến"
$ws1.Range("E16").Value = "Chiến"

# --- Sheet "BigNum" (sheet2): keep the class-name cell as "BigNum" ---
$ws2 = $wb.Worksheets.Item("BigNum")
$ws2.Range("A2").Value = "BigNum"

# Update the saved UI selection state on each sheet. Do the BigNum sheet
# first, then finish on BigInt so it stays the active tab (as in the file).
$ws2.Activate()
[void]$ws2.Range("F6").Select()

$ws1.Activate()
[void]$ws1.Range("E16").Select()
